$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 6.189590430959694

$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 117.745847958593
$ws.Range("D3").Value = 19477208507.93593
$ws.Range("E3").Value = 2195978.878461985
$ws.Range("G3").Value = 19479404607.84708
